$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Dr. Hend Mahmoud, Dr. Alshimaa Atef, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad, Dr. Heba Mahmoud Ali, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Shimaa Ahmad Mekki"
$ws.Range("G3").Value = "Dr. Gehan Adel, Dr. Manar Montaser, Administrator, Dr. Alshimaa Atef"
$ws.Range("G4").Value = "Dr. Nourhan Mahmoud, Dr. Hanan Ragab, Dr. Asmaa Reda, Dr. Majorelle Magdy, Dr. Heba Mahmoud Ali, Dr. Shimaa Ahmad Mekki, Dr. Menna tuâ€™Allah Medhat"
$ws.Range("G6").Value = "Dr. Sara Nabil, Dr. Safa Hany"
$ws.Range("G7").Value = "Dr. Amal Awwad, Dr. Safa Hany"
$ws.Range("G9").Value = "Dr. Marina Youhanna, Dr. Madeha Saeed, Dr. Eman M. Abo-Sakaya, Dr. Yasmeena Fattoh"
$ws.Range("G10").Value = "Dr. Amira Ibrahim, Dr. Basma Hamed"
$ws.Range("G12").Value = "Dr. Mona Ibrahim Hussein, Dr. Dalia Tarek Elwan, Dr. Heba Al-Sayed Mohammad"
$ws.Range("G19").Value = "Dr. Yasmin, Dr. Eman Samir Gabry, Dr. Neveen Nashaat, Dr. Wafaa Ebida, Dr. Marina Sorial"
$ws.Range("G20").Value = "Dr. Yasmin, Dr. Neveen Nashaat, Dr. Nardine, Dr. Remon, Dr. Monica, Dr. Marina Atef, Dr. Marina Sorial"
$ws.Range("G21").Value = "Dr. Hend Mahmoud, Dr. Alshimaa Atef, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad, Dr. Heba Mahmoud Ali, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Shimaa Ahmad Mekki"
$ws.Range("G22").Value = "Dr. Gehan Adel, Dr. Manar Montaser, Administrator, Dr. Alshimaa Atef"
$ws.Range("G23").Value = "Dr. Nourhan Mahmoud, Dr. Hanan Ragab, Dr. Asmaa Reda, Dr. Majorelle Magdy, Dr. Heba Mahmoud Ali, Dr. Shimaa Ahmad Mekki, Dr. Menna tuâ€™Allah Medhat"
$ws.Range("G24").Value = "Dr. Nada Mohammad, Dr. Abeer Ragab, Dr. Lamiaa Ossama, Dr. Fatma Elhady, Dr. Amera Ahmad Saad"
$ws.Range("G25").Value = "Dr. Yasmin Tarek, Dr. Nourhan Mohammad"
$ws.Range("G26").Value = "Dr. Amal Awwad, Dr. Safa Hany"
$ws.Range("G28").Value = "Dr. Marwa Mustafa, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya, Dr. Dina Adel, Dr. Sarah Abdelmohsen, Dr. Esraa Mostafa, Dr. Madeha Saeed, Dr. Basma Hamed, Dr. Nourhan Osama, Dr. Arwa Al-Sayed"
$ws.Range("G29").Value = "Dr. Amira Ibrahim, Dr. Yasmeena Fattoh, Dr. Esraa Mostafa"
$ws.Range("G31").Value = "Dr. Mona Ibrahim Hussein, Dr. Dalia Tarek Elwan, Dr. Heba Al-Sayed Mohammad"
$ws.Range("G38").Value = "Dr. Yasmin, Dr. Neveen Nashaat, Dr. Nardine, Dr. Remon, Dr. Monica, Dr. Marina Atef, Dr. Marina Sorial"
$ws.Range("G39").Value = "Dr. Yasmin, Dr. Neveen Nashaat, Dr. Nardine, Dr. Remon, Dr. Monica, Dr. Marina Atef, Dr. Marina Sorial"
$ws.Range("G40").Value = "Dr. Hend Mahmoud, Dr. Alshimaa Atef, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad, Dr. Heba Mahmoud Ali, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Shimaa Ahmad Mekki"
$ws.Range("G41").Value = "Dr. Hanan Ragab, Dr. Hend Mahmoud, Dr. Alshimaa Atef, Dr. Amira Sobhy, Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki"
$ws.Range("G42").Value = "Dr. Eman Tantawi, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad, Dr. Shimaa Ahmad Mekki, Dr. Menna tuâ€™Allah Medhat"
$ws.Range("G43").Value = "Dr. Nada Mohammad, Dr. Menna tu'Alllah Mohammad, Dr. Abeer Ragab, Dr. Lamiaa Ossama, Dr. Fatma Elhady, Dr. Amera Ahmad Saad, Dr. Kerelos Zareef"
$ws.Range("G44").Value = "Dr. Sara Nabil, Dr. Safa Hany"
$ws.Range("G45").Value = "Dr. Amal Awwad, Dr. Safa Hany"
$ws.Range("G47").Value = "Dr. Merna Said, Dr. Esraa Mostafa, Dr. Amira Ibrahim, Dr. Nourhan Osama, Dr. Maryam Ahmad, Dr. Arwa Al-Sayed"
$ws.Range("G48").Value = "Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya, Dr. Amany Raafat, Dr. Merna Said, Dr. Sarah Abdelmohsen, Dr. Maryam Ahmad"
$ws.Range("G50").Value = "Dr. Mona Ibrahim Hussein, Dr. Dalia Tarek Elwan, Dr. Heba Al-Sayed Mohammad"
$ws.Range("G57").Value = "Dr. Yasmin, Dr. Neveen Nashaat, Dr. Nardine, Dr. Remon, Dr. Monica, Dr. Marina Atef, Dr. Marina Sorial"
$ws.Range("G58").Value = "Dr. Yasmin, Dr. Neveen Nashaat, Dr. Nardine, Dr. Remon, Dr. Monica, Dr. Marina Atef, Dr. Marina Sorial"
$ws.Range("G59").Value = "Dr. Nesma, Dr. Nourhan Mahmoud, Dr. Asmaa Reda, Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad, Dr. Amira Sobhy, Dr. Mohammad El-Tanany"
$ws.Range("G60").Value = "Dr. Hanan Ragab, Dr. Hend Mahmoud, Dr. Alshimaa Atef, Dr. Amira Sobhy, Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki"
$ws.Range("G61").Value = "Dr. Nahla Nagiub, Dr. Asmaa Reda, Dr. Majorelle Magdy, Dr. Amira Sobhy, Dr. Shimaa Ahmad Mekki"
$ws.Range("G63").Value = "Dr. Safa Hany, Dr. Amal Awwad, Dr. Aya Saeed"
$ws.Range("G66").Value = "Dr. Marina Youhanna, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya, Dr. Dina Adel, Dr. Madeha Saeed, Dr. Amira Ibrahim"
$ws.Range("G67").Value = "Dr. Amira Ibrahim, Dr. Yasmeena Fattoh, Dr. Esraa Mostafa"
$ws.Range("G76").Value = "Dr. Yasmin, Dr. Eman Samir Gabry, Dr. Neveen Nashaat, Dr. Wafaa Ebida, Dr. Marina Sorial"
$ws.Range("G77").Value = "Dr. Yasmin, Dr. Neveen Nashaat, Dr. Nardine, Dr. Remon, Dr. Monica, Dr. Marina Atef, Dr. Marina Sorial"
$ws.Range("G78").Value = "Dr. Nesma, Dr. Nourhan Mahmoud, Dr. Asmaa Reda, Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad, Dr. Amira Sobhy, Dr. Mohammad El-Tanany"
$ws.Range("G79").Value = "Dr. Hanan Ragab, Dr. Hend Mahmoud, Dr. Alshimaa Atef, Dr. Amira Sobhy, Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki"
$ws.Range("G80").Value = "Dr. Nahla Nagiub, Dr. Asmaa Reda, Dr. Majorelle Magdy, Dr. Amira Sobhy, Dr. Shimaa Ahmad Mekki"
$ws.Range("G81").Value = "Dr. Nada Mohammad, Dr. Abeer Ragab, Dr. Lamiaa Ossama, Dr. Fatma Elhady, Dr. Amera Ahmad Saad"
$ws.Range("G82").Value = "Dr. Yasmin Tarek, Dr. Nourhan Mohammad"
$ws.Range("G83").Value = "Dr. Aya Saeed, Dr. Amal Awwad, Dr. Safa Hany"
$ws.Range("G85").Value = "Dr. Marina Youhanna, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya, Dr. Dina Adel, Dr. Madeha Saeed, Dr. Amira Ibrahim"
$ws.Range("G86").Value = "Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya, Dr. Amany Raafat, Dr. Merna Said, Dr. Sarah Abdelmohsen, Dr. Maryam Ahmad"
$ws.Range("G88").Value = "Dr. Mona Ibrahim Hussein, Dr. Dalia Tarek Elwan, Dr. Heba Al-Sayed Mohammad"
$ws.Range("G95").Value = "Dr. Yasmin, Dr. Eman Samir Gabry, Dr. Neveen Nashaat, Dr. Wafaa Ebida, Dr. Marina Sorial"
$ws.Range("G96").Value = "Dr. Yasmin, Dr. Neveen Nashaat, Dr. Nardine, Dr. Remon, Dr. Monica, Dr. Marina Atef, Dr. Marina Sorial"
$ws.Range("G97").Value = "Dr. Nesma, Dr. Nourhan Mahmoud, Dr. Asmaa Reda, Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad, Dr. Amira Sobhy, Dr. Mohammad El-Tanany"
$ws.Range("G98").Value = "Dr. Hanan Ragab, Dr. Hend Mahmoud, Dr. Alshimaa Atef, Dr. Amira Sobhy, Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki"
$ws.Range("G99").Value = "Dr. Eman Tantawi, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad, Dr. Shimaa Ahmad Mekki, Dr. Menna tuâ€™Allah Medhat"
$ws.Range("G100").Value = "Dr. Nada Mohammad, Dr. Menna tu'Alllah Mohammad, Dr. Abeer Ragab, Dr. Lamiaa Ossama, Dr. Fatma Elhady, Dr. Amera Ahmad Saad, Dr. Kerelos Zareef"
$ws.Range("G101").Value = "Dr. Safa Hany, Dr. Amal Awwad, Dr. Aya Saeed"
$ws.Range("G102").Value = "Dr. Amal Awwad, Dr. Safa Hany"
$ws.Range("G104").Value = "Dr. Merna Said, Dr. Esraa Mostafa, Dr. Amira Ibrahim, Dr. Nourhan Osama, Dr. Maryam Ahmad, Dr. Arwa Al-Sayed"
$ws.Range("G105").Value = "Dr. Amira Ibrahim, Dr. Basma Hamed"
$ws.Range("G115").Value = "Dr. Yasmin, Dr. Neveen Nashaat, Dr. Nardine, Dr. Remon, Dr. Monica, Dr. Marina Atef, Dr. Marina Sorial"
